$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("B2").Value = "Não será mexida"
$ws.Range("B3").Value = "Resolvido"
$ws.Range("C3").Value = "Victor"
$ws.Range("B23").Value = "-"
$ws.Range("C23").Value = "-"

$ws.Activate()
$ws.Range("C24").Select()
